$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.327452000000001
$ws.Range("H2").Value = 21.982356
$ws.Range("I2").Value = 0.2396705957528817
$ws.Range("J2").Value = 0.2396705957528817
$ws.Range("M2").Value = 125.901487
$ws.Range("N2").Value = 377.704461
$ws.Range("O2").Value = 0.8798726812012091
$ws.Range("P2").Value = 0.879872681201209
$ws.Range("Q2").Value = 922.5371027211243
$ws.Range("R2").Value = 8302.833924490118
$ws.Range("S2").Value = 0.2108796096901792
$ws.Range("T2").Value = 0.2108796096901791
$ws.Range("G3").Value = 7.327452000000001
$ws.Range("H3").Value = 21.982356
$ws.Range("I3").Value = 0.2396705957528817
$ws.Range("J3").Value = 0.2396705957528817
$ws.Range("O3").Value = 0.006089432091685741
$ws.Range("P3").Value = 0.006089432091685741
$ws.Range("Q3").Value = 6.384704468164001
$ws.Range("R3").Value = 57.46234021347601
$ws.Range("S3").Value = 0.001459457817211038
$ws.Range("T3").Value = 0.001459457817211038
$ws.Range("G4").Value = 7.327452000000001
$ws.Range("H4").Value = 21.982356
$ws.Range("I4").Value = 0.2396705957528817
$ws.Range("J4").Value = 0.2396705957528817
$ws.Range("M4").Value = 2.781571666666667
$ws.Range("N4").Value = 8.344715000000001
$ws.Range("O4").Value = 0.01943923760251788
$ws.Range("P4").Value = 0.01943923760251788
$ws.Range("Q4").Value = 20.38183287206
$ws.Range("R4").Value = 183.43649584854
$ws.Range("S4").Value = 0.00465901365717728
$ws.Range("T4").Value = 0.004659013657177279
$ws.Range("G5").Value = 7.327452000000001
$ws.Range("H5").Value = 21.982356
$ws.Range("I5").Value = 0.2396705957528817
$ws.Range("J5").Value = 0.2396705957528817
$ws.Range("M5").Value = 13.536175
$ws.Range("N5").Value = 40.608525
$ws.Range("O5").Value = 0.09459864910458742
$ws.Range("P5").Value = 0.09459864910458742
$ws.Range("Q5").Value = 99.18567257610002
$ws.Range("R5").Value = 892.6710531849001
$ws.Range("S5").Value = 0.02267251458831428
$ws.Range("T5").Value = 0.02267251458831427
$ws.Range("I6").Value = 0.3687153231745697
$ws.Range("J6").Value = 0.3687153231745697
$ws.Range("M6").Value = 125.901487
$ws.Range("N6").Value = 377.704461
$ws.Range("O6").Value = 0.8798726812012091
$ws.Range("P6").Value = 0.879872681201209
$ws.Range("Q6").Value = 1419.254476761406
$ws.Range("R6").Value = 12773.29029085265
$ws.Range("S6").Value = 0.324422540001579
$ws.Range("T6").Value = 0.3244225400015789
$ws.Range("I7").Value = 0.3687153231745697
$ws.Range("J7").Value = 0.3687153231745697
$ws.Range("O7").Value = 0.006089432091685741
$ws.Range("P7").Value = 0.006089432091685741
$ws.Range("S7").Value = 0.002245266921635504
$ws.Range("T7").Value = 0.002245266921635504
$ws.Range("I8").Value = 0.3687153231745697
$ws.Range("J8").Value = 0.3687153231745697
$ws.Range("M8").Value = 2.781571666666667
$ws.Range("N8").Value = 8.344715000000001
$ws.Range("O8").Value = 0.01943923760251788
$ws.Range("P8").Value = 0.01943923760251788
$ws.Range("Q8").Value = 31.35592862655666
$ws.Range("R8").Value = 282.20335763901
$ws.Range("S8").Value = 0.007167544774879627
$ws.Range("T8").Value = 0.007167544774879627
$ws.Range("I9").Value = 0.3687153231745697
$ws.Range("J9").Value = 0.3687153231745697
$ws.Range("M9").Value = 13.536175
$ws.Range("N9").Value = 40.608525
$ws.Range("O9").Value = 0.09459864910458742
$ws.Range("P9").Value = 0.09459864910458742
$ws.Range("Q9").Value = 152.58975429715
$ws.Range("R9").Value = 1373.30778867435
$ws.Range("S9").Value = 0.03487997147647567
$ws.Range("T9").Value = 0.03487997147647567
$ws.Range("G10").Value = 6.386255666666667
$ws.Range("H10").Value = 19.158767
$ws.Range("I10").Value = 0.2088853943035337
$ws.Range("J10").Value = 0.2088853943035337
$ws.Range("M10").Value = 125.901487
$ws.Range("N10").Value = 377.704461
$ws.Range("O10").Value = 0.8798726812012091
$ws.Range("P10").Value = 0.879872681201209
$ws.Range("Q10").Value = 804.0390847955098
$ws.Range("R10").Value = 7236.351763159588
$ws.Range("S10").Value = 0.183792551949622
$ws.Range("T10").Value = 0.1837925519496219
$ws.Range("G11").Value = 6.386255666666667
$ws.Range("H11").Value = 19.158767
$ws.Range("I11").Value = 0.2088853943035337
$ws.Range("J11").Value = 0.2088853943035337
$ws.Range("O11").Value = 0.006089432091685741
$ws.Range("P11").Value = 0.006089432091685741
$ws.Range("Q11").Value = 5.564602141345222
$ws.Range("R11").Value = 50.08141927210701
$ws.Range("S11").Value = 0.001271993423556368
$ws.Range("T11").Value = 0.001271993423556368
$ws.Range("G12").Value = 6.386255666666667
$ws.Range("H12").Value = 19.158767
$ws.Range("I12").Value = 0.2088853943035337
$ws.Range("J12").Value = 0.2088853943035337
$ws.Range("M12").Value = 2.781571666666667
$ws.Range("N12").Value = 8.344715000000001
$ws.Range("O12").Value = 0.01943923760251788
$ws.Range("P12").Value = 0.01943923760251788
$ws.Range("Q12").Value = 17.76382781848945
$ws.Range("R12").Value = 159.874450366405
$ws.Range("S12").Value = 0.004060572811562026
$ws.Range("T12").Value = 0.004060572811562025
$ws.Range("G13").Value = 6.386255666666667
$ws.Range("H13").Value = 19.158767
$ws.Range("I13").Value = 0.2088853943035337
$ws.Range("J13").Value = 0.2088853943035337
$ws.Range("M13").Value = 13.536175
$ws.Range("N13").Value = 40.608525
$ws.Range("O13").Value = 0.09459864910458742
$ws.Range("P13").Value = 0.09459864910458742
$ws.Range("Q13").Value = 86.44547429874167
$ws.Range("R13").Value = 778.0092686886751
$ws.Range("S13").Value = 0.01976027611879337
$ws.Range("T13").Value = 0.01976027611879336
$ws.Range("G14").Value = 5.586566333333333
$ws.Range("H14").Value = 16.759699
$ws.Range("I14").Value = 0.1827286867690149
$ws.Range("J14").Value = 0.1827286867690149
$ws.Range("M14").Value = 125.901487
$ws.Range("N14").Value = 377.704461
$ws.Range("O14").Value = 0.8798726812012091
$ws.Range("P14").Value = 0.879872681201209
$ws.Range("Q14").Value = 703.3570085908044
$ws.Range("R14").Value = 6330.213077317238
$ws.Range("S14").Value = 0.160777979559829
$ws.Range("T14").Value = 0.160777979559829
$ws.Range("G15").Value = 5.586566333333333
$ws.Range("H15").Value = 16.759699
$ws.Range("I15").Value = 0.1827286867690149
$ws.Range("J15").Value = 0.1827286867690149
$ws.Range("O15").Value = 0.006089432091685741
$ws.Range("P15").Value = 0.006089432091685741
$ws.Range("Q15").Value = 4.867800571075445
$ws.Range("R15").Value = 43.810205139679
$ws.Range("S15").Value = 0.001112713929282831
$ws.Range("T15").Value = 0.001112713929282831
$ws.Range("G16").Value = 5.586566333333333
$ws.Range("H16").Value = 16.759699
$ws.Range("I16").Value = 0.1827286867690149
$ws.Range("J16").Value = 0.1827286867690149
$ws.Range("M16").Value = 2.781571666666667
$ws.Range("N16").Value = 8.344715000000001
$ws.Range("O16").Value = 0.01943923760251788
$ws.Range("P16").Value = 0.01943923760251788
$ws.Range("Q16").Value = 15.53943462675389
$ws.Range("R16").Value = 139.854911640785
$ws.Range("S16").Value = 0.003552106358898945
$ws.Range("T16").Value = 0.003552106358898945
$ws.Range("G17").Value = 5.586566333333333
$ws.Range("H17").Value = 16.759699
$ws.Range("I17").Value = 0.1827286867690149
$ws.Range("J17").Value = 0.1827286867690149
$ws.Range("M17").Value = 13.536175
$ws.Range("N17").Value = 40.608525
$ws.Range("O17").Value = 0.09459864910458742
$ws.Range("P17").Value = 0.09459864910458742
$ws.Range("Q17").Value = 75.62073953710832
$ws.Range("R17").Value = 680.5866558339749
$ws.Range("S17").Value = 0.0172858869210041
$ws.Range("T17").Value = 0.0172858869210041

Write-Output "Applied 174 cell updates"
